# "Add new skills to database directly"
#
# The "U-Net" entry (Skill = "U-Net", Domain = "Data Science") is removed
# from the "Skills" table. Deleting the whole worksheet row shifts every
# row below it up by one, which is exactly what turns the 93-row table
# (A1:B93) into a 92-row table (A1:B92), drops "U-Net" from the shared
# strings, and leaves the selection sitting on the row that used to hold
# U-Net (now occupied by the former row 24, "MobileNet").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skills")

# Row 23 is "U-Net" (Skill column A23, Domain column B23).
$ws.Rows.Item(23).Select()
$ws.Rows.Item(23).Delete()
